# fix: alterar python version para 3.11.5
# Update absenteeism data rows 2-11 on the active sheet with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 41943
$ws.Cells.Item(2, 2).Value = "Sr. Yuri Cavalcanti"
$ws.Cells.Item(2, 3).Value = "TI"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 45102
$ws.Cells.Item(2, 7).Value = 3615.67

# Row 3
$ws.Cells.Item(3, 1).Value = 10224
$ws.Cells.Item(3, 2).Value = "Laura Teixeira"
$ws.Cells.Item(3, 4).Value = "Consulta médica"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45082
$ws.Cells.Item(3, 7).Value = 2738.03

# Row 4
$ws.Cells.Item(4, 1).Value = 88657
$ws.Cells.Item(4, 2).Value = "Vicente Jesus"
$ws.Cells.Item(4, 3).Value = "Recursos Humanos"
$ws.Cells.Item(4, 4).Value = "Outros"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 45086
$ws.Cells.Item(4, 7).Value = 11212.27

# Row 5
$ws.Cells.Item(5, 1).Value = 78671
$ws.Cells.Item(5, 2).Value = "Diego Moreira"
$ws.Cells.Item(5, 3).Value = "Marketing"
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = 45106
$ws.Cells.Item(5, 7).Value = 6762.71

# Row 6
$ws.Cells.Item(6, 1).Value = 94072
$ws.Cells.Item(6, 2).Value = "Ana Luiza Monteiro"
$ws.Cells.Item(6, 3).Value = "Financeiro"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 45101
$ws.Cells.Item(6, 7).Value = 3455.53

# Row 7
$ws.Cells.Item(7, 1).Value = 877
$ws.Cells.Item(7, 2).Value = "Pietra Monteiro"
$ws.Cells.Item(7, 4).Value = "Viagem de negócios"
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 45080
$ws.Cells.Item(7, 7).Value = 7473.84

# Row 8
$ws.Cells.Item(8, 1).Value = 17350
$ws.Cells.Item(8, 2).Value = "Milena Silveira"
$ws.Cells.Item(8, 3).Value = "TI"
$ws.Cells.Item(8, 4).Value = "Problemas pessoais"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45099
$ws.Cells.Item(8, 7).Value = 9037.629999999999

# Row 9
$ws.Cells.Item(9, 1).Value = 80753
$ws.Cells.Item(9, 2).Value = "Bianca Correia"
$ws.Cells.Item(9, 3).Value = "Vendas"
$ws.Cells.Item(9, 4).Value = "Problemas pessoais"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 45094
$ws.Cells.Item(9, 7).Value = 4931.86

# Row 10
$ws.Cells.Item(10, 1).Value = 91190
$ws.Cells.Item(10, 2).Value = "Sra. Sophie Moraes"
$ws.Cells.Item(10, 4).Value = "Consulta médica"
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 45089
$ws.Cells.Item(10, 7).Value = 6826.45

# Row 11
$ws.Cells.Item(11, 1).Value = 4202
$ws.Cells.Item(11, 2).Value = "Dr. Luiz Otávio Teixeira"
$ws.Cells.Item(11, 3).Value = "Marketing"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 45104
$ws.Cells.Item(11, 7).Value = 2609.09
